$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 1: "Министерство образования Московской области"
#           -> "Министерство образования " | "ХХХ"   (2 runs)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range

# Find/Replace on this engine silently re-merges the replacement text back
# into a neighbouring, identically formatted run. Work around it by giving
# the replacement a throw-away direct-formatting difference (Bold) so it
# is forced to stay in its own run, then clear that throw-away formatting
# again afterwards (Bold is the one property this engine fully removes
# from <w:rPr> again instead of emitting an explicit "off" toggle).
$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Replacement.Font.Bold = $true
$p1.Find.Execute("Московской области", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "ХХХ", 2)

# $p1 auto-collapsed to the freshly inserted "ХХХ" text - clear the
# throw-away Bold directly on it.
$p1.Font.Bold = $false

# ---------------------------------------------------------------------
# Paragraph 3: "Московской области «" | "ХХХХ" | "»"
#           -> "ХХХ" | " «" | "ХХХХ" | "»"            (4 runs)
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$S = $p3.Start

$oldText = "Московской области"
$newText = "ХХХ"
$run1Len = $oldText.Length + 2          # "Московской области" + " «"
$run1Range = $d.Range($S, $S + $run1Len)

$d.Content.Find.ClearFormatting()
$d.Content.Find.Replacement.ClearFormatting()
$d.Content.Find.Replacement.Font.Bold = $true
$run1Range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)

# $run1Range auto-collapsed to the freshly inserted "ХХХ" text.
$run1Range.Font.Bold = $false
$newTextEnd = $S + $newText.Length

# The untouched tail (" «" + "ХХХХ" + "»") got coalesced into a single
# run by the replace above. Re-split it back into its three original
# pieces using the same throw-away-Bold trick, one boundary at a time.
$segQuoteStart = $newTextEnd
$segQuoteEnd = $segQuoteStart + 2        # " «"
$segXEnd = $segQuoteEnd + 4              # "ХХХХ"
$segCloseEnd = $segXEnd + 1              # "»"

$segQuote = $d.Range($segQuoteStart, $segQuoteEnd)
$segQuote.Font.Bold = $true
$segClose = $d.Range($segXEnd, $segCloseEnd)
$segClose.Font.Bold = $true

$segQuote2 = $d.Range($segQuoteStart, $segQuoteEnd)
$segQuote2.Font.Bold = $false
$segClose2 = $d.Range($segXEnd, $segCloseEnd)
$segClose2.Font.Bold = $false

$d.Save()
